# Updates cryptos list prices / 1h-volume figures (and fixes the row 21/22
# ordering of ShibaInu vs InternetComputer(DFINITY)) to match the latest
# GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.199.15"
$ws.Range("E2").Value = "  -0.77%  "

# Row 3
$ws.Range("D3").Value = "3.064.14"
$ws.Range("E3").Value = "  +1.33%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "387.72"
$ws.Range("E5").Value = "  +2.16%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.11"
$ws.Range("E6").Value = "  -0.12%  "

# Row 7
$ws.Range("E7").Value = "  -1.69%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("E9").Value = "  -1.63%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.77"
$ws.Range("E10").Value = "  +0.18%  "

# Row 11
$ws.Range("E11").Value = "  +0.25%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0848"
$ws.Range("E12").Value = "  -1.33%  "

# Row 13
$ws.Range("D13").Value = "3.542.22"
$ws.Range("E13").Value = "  +1.04%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.28"
$ws.Range("E14").Value = "  -1.05%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.69"
$ws.Range("E15").Value = "  -0.45%  "

# Row 16
$ws.Range("D16").Value = "3.058.23"
$ws.Range("E16").Value = "  +0.78%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.991"
$ws.Range("E17").Value = "  +1.91%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.65"
$ws.Range("E18").Value = "  +0.35%  "

# Row 19
$ws.Range("D19").Value = "51.198.30"
$ws.Range("E19").Value = "  -0.84%  "

# Row 20
$ws.Range("E20").Value = "  +3.26%  "

# Row 21
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.28"
$ws.Range("E21").Value = "  -1.01%  "

# Row 22
$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").Value = "0.0₃0955"
$ws.Range("E22").Value = "  -0.68%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.66"
$ws.Range("E23").Value = "  -0.39%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.50"
$ws.Range("E24").Value = "  -0.96%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.11"
$ws.Range("E25").Value = "  -0.93%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.91"
$ws.Range("E26").Value = "  -4.99%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.09"
$ws.Range("E27").Value = "  +3.62%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.25"
$ws.Range("E28").Value = "  -3.74%  "

# Row 29
$ws.Range("E29").Value = "  +0.02%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.163"
$ws.Range("E30").Value = "  -5.14%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.105"
$ws.Range("E31").Value = "  -3.20%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.41"
$ws.Range("E32").Value = "  +1.66%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.66"
$ws.Range("E33").Value = "  +5.18%  "

# Row 34
$ws.Range("E34").Value = "  +4.89%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.08"
$ws.Range("E35").Value = "  +1.83%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.05"
$ws.Range("E36").Value = "  -0.97%  "

# Row 37
$ws.Range("E37").Value = "  -0.14%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.38"
$ws.Range("E38").Value = "  +2.66%  "

# Row 39
$ws.Range("E39").Value = "  -0.27%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "130.55"
$ws.Range("E40").Value = "  +2.28%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.83"
$ws.Range("E41").Value = "  -1.18%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.52"
$ws.Range("E42").Value = "  -2.73%  "

# Row 43
$ws.Range("E43").Value = "  -0.58%  "

# Row 44
$ws.Range("E44").Value = "  +0.19%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.48"
$ws.Range("E45").Value = "  -1.95%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.65"
$ws.Range("E46").Value = "  +0.22%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.52"
$ws.Range("E47").Value = "  +3.69%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.06"
$ws.Range("E48").Value = "  -0.12%  "

# Row 49
$ws.Range("D49").Value = "2.058.44"
$ws.Range("E49").Value = "  +1.92%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0326"
$ws.Range("E50").Value = "  +3.30%  "

# Row 51
$ws.Range("E51").Value = "  +14.28%  "
